{"js": "// Apply the wording edits to the \"\u0410\u043d\u043a\u0435\u0442\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u0435\" questionnaire document.\n\nasync function replaceOnce(body, searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. Purpose paragraph: narrow the product's target audience wording.\nawait replaceOnce(\n  body,\n  \"\u0432\u044b\u044f\u0441\u043d\u0438\u0442\u044c, \u0442\u0440\u0435\u0431\u0443\u0435\u0442\u0441\u044f \u043b\u0438 \u043f\u0440\u0438\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0434\u043b\u044f \u043b\u044e\u0434\u0435\u0439 \u0441 \u0434\u0435\u0444\u0435\u043a\u0442\u0430\u043c\u0438 \u0437\u0440\u0435\u043d\u0438\u044f.\",\n  \"\u0432\u044b\u044f\u0441\u043d\u0438\u0442\u044c, \u0442\u0440\u0435\u0431\u0443\u0435\u0442\u0441\u044f \u043b\u0438 \u043f\u0440\u0438\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0434\u043b\u044f \u0441\u0432\u044f\u0437\u0438 \u043b\u044e\u0434\u0435\u0439 \u0441 \u043f\u043e\u043c\u043e\u0449\u044c\u044e \u044e\u043c\u043e\u0440\u0430\\\\\u043c\u0435\u043c\u043e\u0432.\"\n);\n\n// 2. Question: vision defects -> communication / shared interests problems.\nawait replaceOnce(\n  body,\n  \"\u0415\u0441\u0442\u044c \u043b\u0438 \u0443 \u0412\u0430\u0441 \u0438\u043b\u0438 \u0443 \u0412\u0430\u0448\u0438\u0445 \u0437\u043d\u0430\u043a\u043e\u043c\u044b\u0445 \u0434\u0435\u0444\u0435\u043a\u0442\u044b \u0437\u0440\u0435\u043d\u0438\u044f?\",\n  \"\u0415\u0441\u0442\u044c \u043b\u0438 \u0443 \u0412\u0430\u0441 \u0438\u043b\u0438 \u0443 \u0412\u0430\u0448\u0438\u0445 \u0437\u043d\u0430\u043a\u043e\u043c\u044b\u0445 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u044b \u0441 \u043a\u043e\u043c\u043c\u0443\u043d\u0438\u043a\u0430\u0446\u0438\u0435\u0439 \u0438 \u043d\u0430\u0445\u043e\u0436\u0434\u0435\u043d\u0438\u0435\u043c \u043e\u0431\u0449\u0438\u0445 \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043e\u0432?\"\n);\n\n// 3. Question: attitude towards people with vision problems -> differing interests.\nawait replaceOnce(\n  body,\n  \"\u041a\u0430\u043a\u043e\u0435 \u0443 \u0412\u0430\u0441 \u043e\u0442\u043d\u043e\u0448\u0435\u043d\u0438\u0435 \u043a \u043b\u044e\u0434\u044f\u043c \u0441\u043e \u0437\u043d\u0430\u0447\u0438\u0442\u0435\u043b\u044c\u043d\u044b\u043c\u0438 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u0430\u043c\u0438 \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c?\",\n  \"\u041a\u0430\u043a\u043e\u0435 \u0443 \u0412\u0430\u0441 \u043e\u0442\u043d\u043e\u0448\u0435\u043d\u0438\u0435 \u043a \u043b\u044e\u0434\u044f\u043c \u043a \u043b\u044e\u0434\u044f\u043c \u0438\u043c\u0435\u044e\u0449\u0438\u043c \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u044b \u043e\u0442\u043b\u0438\u0447\u0430\u044e\u0449\u0438\u0445\u0441\u044f \u043e\u0442 \u0412\u0430\u0448\u0438\u0445 \u0441\u043e\u0431\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0445?\"\n);\n\n// 4. Question: messenger purpose wording.\nawait replaceOnce(\n  body,\n  \"\u041d\u0443\u0436\u0435\u043d \u043b\u0438 \u043c\u0438\u0440\u0443 \u043c\u0435\u0441\u0441\u0435\u043d\u0434\u0436\u0435\u0440, \u043a\u043e\u0442\u043e\u0440\u044b\u0439 \u0431\u044b \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u043b \u043e\u0431\u0449\u0430\u0442\u044c\u0441\u044f \u0438 \u0434\u0435\u043b\u0438\u0442\u044c\u0441\u044f \u043c\u044b\u0441\u043b\u044f\u043c\u0438 \u043a\u0430\u043a \u043e\u0431\u044b\u0447\u043d\u044b\u043c \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0418\u043d\u0442\u0435\u0440\u043d\u0435\u0442\u0430, \u0442\u0430\u043a \u0438 \u043b\u044e\u0434\u044f\u043c \u0441 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u0430\u043c\u0438 \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c?\",\n  \"\u041d\u0443\u0436\u0435\u043d \u043b\u0438 \u043c\u0438\u0440\u0443 \u043c\u0435\u0441\u0441\u0435\u043d\u0434\u0436\u0435\u0440\\\\\u0441\u043e\u0446\u0435\u0441\u0435\u0442\u044c, \u043a\u043e\u0442\u043e\u0440\u044b\u0439 \u0431\u044b \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u043b \u043e\u0431\u0449\u0430\u0442\u044c\u0441\u044f \u0438 \u0434\u0435\u043b\u0438\u0442\u044c\u0441\u044f \u043c\u044b\u0441\u043b\u044f\u043c\u0438 \u043f\u043e\u0441\u0440\u0435\u0434\u0441\u0442\u0432\u043e\u043c \u044e\u043c\u043e\u0440\u0430 \u0438 \u043c\u0435\u043c\u043e\u0432\\\\\u0440\u043e\u0444\u043b\u043e\u0432?\"\n);\n\n// 5. Question: wish to keep in touch with people who share a sense of humor.\nawait replaceOnce(\n  body,\n  \"\u0411\u044b\u043b\u043e \u0431\u044b \u0443 \u0432\u0430\u0441 \u0436\u0435\u043b\u0430\u043d\u0438\u0435 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0438\u0432\u0430\u0442\u044c \u043e\u0431\u0449\u0435\u043d\u0438\u0435 \u0438\u043b\u0438 \u0437\u043d\u0430\u043a\u043e\u043c\u0438\u0442\u044c\u0441\u044f \u0441 \u043b\u044e\u0434\u044c\u043c\u0438, \u0443 \u043a\u043e\u0442\u043e\u0440\u044b\u0445 \u0435\u0441\u0442\u044c \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u044b \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c, \u0447\u0435\u0440\u0435\u0437 \u0438\u043d\u0442\u0435\u0440\u043d\u0435\u0442\",\n  \"\u0411\u044b\u043b\u043e \u0431\u044b \u0443 \u0432\u0430\u0441 \u0436\u0435\u043b\u0430\u043d\u0438\u0435 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0438\u0432\u0430\u0442\u044c \u043e\u0431\u0449\u0435\u043d\u0438\u0435 \u0438\u043b\u0438 \u0437\u043d\u0430\u043a\u043e\u043c\u0438\u0442\u044c\u0441\u044f \u0441 \u043b\u044e\u0434\u044c\u043c\u0438 \u0443 \u043a\u043e\u0442\u043e\u0440\u044b\u0445 \u0447\u0443\u0432\u0441\u0442\u0432\u043e \u044e\u043c\u043e\u0440\u0430 \u0441\u0445\u043e\u0436\u0435 \u0441 \u0432\u0430\u0448\u0438\u043c\"\n);\n\n// 6. Question: feed with humorous posts instead of voice posts.\nawait replaceOnce(\n  body,\n  \"\u0411\u044b\u043b\u0430 \u0431\u044b \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u0430 \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0441\u043e\u0446\u0441\u0435\u0442\u044c, \u0432 \u043b\u0435\u043d\u0442\u0435 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u043f\u0440\u0435\u043e\u0431\u043b\u0430\u0434\u0430\u043b\u0438 \u0431\u044b \u0433\u043e\u043b\u043e\u0441\u043e\u0432\u044b\u0435 \u043f\u043e\u0441\u0442\u044b\",\n  \"\u0411\u044b\u043b\u0430 \u0431\u044b \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u0430 \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0441\u043e\u0446\u0441\u0435\u0442\u044c, \u0432 \u043b\u0435\u043d\u0442\u0435 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u043f\u0440\u0435\u043e\u0431\u043b\u0430\u0434\u0430\u043b\u0438 \u0431\u044b \u044e\u043c\u043e\u0440\u0438\u0441\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u043f\u043e\u0441\u0442\u044b\"\n);\n\n// 7. Add a \"//TODO\" placeholder paragraph right after the page-break\n// paragraph that precedes the results table.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet breakParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  // The paragraph holding the page break renders as a lone form-feed (\\f)\n  // character and directly follows the \"\u0421\u043f\u0430\u0441\u0438\u0431\u043e \u0437\u0430 \u043f\u043e\u043c\u043e\u0449\u044c!\" paragraph.\n  if (p.text.indexOf(\"\\f\") !== -1 && i > 0 && paragraphs.items[i - 1].text.indexOf(\"\u0421\u043f\u0430\u0441\u0438\u0431\u043e \u0437\u0430 \u043f\u043e\u043c\u043e\u0449\u044c\") !== -1) {\n    breakParagraph = p;\n    break;\n  }\n}\n\nif (!breakParagraph) {\n  throw new Error(\"Could not locate the page-break paragraph before the table.\");\n}\n\nconst todoParagraph = breakParagraph.insertParagraph(\"//TODO\", Word.InsertLocation.after);\ntodoParagraph.font.bold = true;\ntodoParagraph.font.italic = true;\ntodoParagraph.font.set({ name: \"Times New Roman\", size: 12 });\nawait context.sync();\n", "ps1": "# Apply the wording edits to the \"\u0410\u043d\u043a\u0435\u0442\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u0435\" questionnaire document.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Search text not found: $findText\"\n    }\n}\n\n# 1. Purpose paragraph: narrow the product's target audience wording.\nReplace-Text \"\u0432\u044b\u044f\u0441\u043d\u0438\u0442\u044c, \u0442\u0440\u0435\u0431\u0443\u0435\u0442\u0441\u044f \u043b\u0438 \u043f\u0440\u0438\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0434\u043b\u044f \u043b\u044e\u0434\u0435\u0439 \u0441 \u0434\u0435\u0444\u0435\u043a\u0442\u0430\u043c\u0438 \u0437\u0440\u0435\u043d\u0438\u044f.\" \"\u0432\u044b\u044f\u0441\u043d\u0438\u0442\u044c, \u0442\u0440\u0435\u0431\u0443\u0435\u0442\u0441\u044f \u043b\u0438 \u043f\u0440\u0438\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u0434\u043b\u044f \u0441\u0432\u044f\u0437\u0438 \u043b\u044e\u0434\u0435\u0439 \u0441 \u043f\u043e\u043c\u043e\u0449\u044c\u044e \u044e\u043c\u043e\u0440\u0430\\\u043c\u0435\u043c\u043e\u0432.\"\n\n# 2. Question: vision defects -> communication / shared interests problems.\nReplace-Text \"\u0415\u0441\u0442\u044c \u043b\u0438 \u0443 \u0412\u0430\u0441 \u0438\u043b\u0438 \u0443 \u0412\u0430\u0448\u0438\u0445 \u0437\u043d\u0430\u043a\u043e\u043c\u044b\u0445 \u0434\u0435\u0444\u0435\u043a\u0442\u044b \u0437\u0440\u0435\u043d\u0438\u044f?\" \"\u0415\u0441\u0442\u044c \u043b\u0438 \u0443 \u0412\u0430\u0441 \u0438\u043b\u0438 \u0443 \u0412\u0430\u0448\u0438\u0445 \u0437\u043d\u0430\u043a\u043e\u043c\u044b\u0445 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u044b \u0441 \u043a\u043e\u043c\u043c\u0443\u043d\u0438\u043a\u0430\u0446\u0438\u0435\u0439 \u0438 \u043d\u0430\u0445\u043e\u0436\u0434\u0435\u043d\u0438\u0435\u043c \u043e\u0431\u0449\u0438\u0445 \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043e\u0432?\"\n\n# 3. Question: attitude towards people with vision problems -> differing interests.\nReplace-Text \"\u041a\u0430\u043a\u043e\u0435 \u0443 \u0412\u0430\u0441 \u043e\u0442\u043d\u043e\u0448\u0435\u043d\u0438\u0435 \u043a \u043b\u044e\u0434\u044f\u043c \u0441\u043e \u0437\u043d\u0430\u0447\u0438\u0442\u0435\u043b\u044c\u043d\u044b\u043c\u0438 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u0430\u043c\u0438 \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c?\" \"\u041a\u0430\u043a\u043e\u0435 \u0443 \u0412\u0430\u0441 \u043e\u0442\u043d\u043e\u0448\u0435\u043d\u0438\u0435 \u043a \u043b\u044e\u0434\u044f\u043c \u043a \u043b\u044e\u0434\u044f\u043c \u0438\u043c\u0435\u044e\u0449\u0438\u043c \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u044b \u043e\u0442\u043b\u0438\u0447\u0430\u044e\u0449\u0438\u0445\u0441\u044f \u043e\u0442 \u0412\u0430\u0448\u0438\u0445 \u0441\u043e\u0431\u0441\u0442\u0432\u0435\u043d\u043d\u044b\u0445?\"\n\n# 4. Question: messenger purpose wording.\nReplace-Text \"\u041d\u0443\u0436\u0435\u043d \u043b\u0438 \u043c\u0438\u0440\u0443 \u043c\u0435\u0441\u0441\u0435\u043d\u0434\u0436\u0435\u0440, \u043a\u043e\u0442\u043e\u0440\u044b\u0439 \u0431\u044b \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u043b \u043e\u0431\u0449\u0430\u0442\u044c\u0441\u044f \u0438 \u0434\u0435\u043b\u0438\u0442\u044c\u0441\u044f \u043c\u044b\u0441\u043b\u044f\u043c\u0438 \u043a\u0430\u043a \u043e\u0431\u044b\u0447\u043d\u044b\u043c \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0418\u043d\u0442\u0435\u0440\u043d\u0435\u0442\u0430, \u0442\u0430\u043a \u0438 \u043b\u044e\u0434\u044f\u043c \u0441 \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u0430\u043c\u0438 \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c?\" \"\u041d\u0443\u0436\u0435\u043d \u043b\u0438 \u043c\u0438\u0440\u0443 \u043c\u0435\u0441\u0441\u0435\u043d\u0434\u0436\u0435\u0440\\\u0441\u043e\u0446\u0435\u0441\u0435\u0442\u044c, \u043a\u043e\u0442\u043e\u0440\u044b\u0439 \u0431\u044b \u043f\u043e\u0437\u0432\u043e\u043b\u044f\u043b \u043e\u0431\u0449\u0430\u0442\u044c\u0441\u044f \u0438 \u0434\u0435\u043b\u0438\u0442\u044c\u0441\u044f \u043c\u044b\u0441\u043b\u044f\u043c\u0438 \u043f\u043e\u0441\u0440\u0435\u0434\u0441\u0442\u0432\u043e\u043c \u044e\u043c\u043e\u0440\u0430 \u0438 \u043c\u0435\u043c\u043e\u0432\\\u0440\u043e\u0444\u043b\u043e\u0432?\"\n\n# 5. Question: wish to keep in touch with people who share a sense of humor.\nReplace-Text \"\u0411\u044b\u043b\u043e \u0431\u044b \u0443 \u0432\u0430\u0441 \u0436\u0435\u043b\u0430\u043d\u0438\u0435 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0438\u0432\u0430\u0442\u044c \u043e\u0431\u0449\u0435\u043d\u0438\u0435 \u0438\u043b\u0438 \u0437\u043d\u0430\u043a\u043e\u043c\u0438\u0442\u044c\u0441\u044f \u0441 \u043b\u044e\u0434\u044c\u043c\u0438, \u0443 \u043a\u043e\u0442\u043e\u0440\u044b\u0445 \u0435\u0441\u0442\u044c \u043f\u0440\u043e\u0431\u043b\u0435\u043c\u044b \u0441\u043e \u0437\u0440\u0435\u043d\u0438\u0435\u043c, \u0447\u0435\u0440\u0435\u0437 \u0438\u043d\u0442\u0435\u0440\u043d\u0435\u0442\" \"\u0411\u044b\u043b\u043e \u0431\u044b \u0443 \u0432\u0430\u0441 \u0436\u0435\u043b\u0430\u043d\u0438\u0435 \u043f\u043e\u0434\u0434\u0435\u0440\u0436\u0438\u0432\u0430\u0442\u044c \u043e\u0431\u0449\u0435\u043d\u0438\u0435 \u0438\u043b\u0438 \u0437\u043d\u0430\u043a\u043e\u043c\u0438\u0442\u044c\u0441\u044f \u0441 \u043b\u044e\u0434\u044c\u043c\u0438 \u0443 \u043a\u043e\u0442\u043e\u0440\u044b\u0445 \u0447\u0443\u0432\u0441\u0442\u0432\u043e \u044e\u043c\u043e\u0440\u0430 \u0441\u0445\u043e\u0436\u0435 \u0441 \u0432\u0430\u0448\u0438\u043c\"\n\n# 6. Question: feed with humorous posts instead of voice posts.\nReplace-Text \"\u0411\u044b\u043b\u0430 \u0431\u044b \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u0430 \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0441\u043e\u0446\u0441\u0435\u0442\u044c, \u0432 \u043b\u0435\u043d\u0442\u0435 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u043f\u0440\u0435\u043e\u0431\u043b\u0430\u0434\u0430\u043b\u0438 \u0431\u044b \u0433\u043e\u043b\u043e\u0441\u043e\u0432\u044b\u0435 \u043f\u043e\u0441\u0442\u044b\" \"\u0411\u044b\u043b\u0430 \u0431\u044b \u0438\u043d\u0442\u0435\u0440\u0435\u0441\u043d\u0430 \u043f\u043e\u043b\u044c\u0437\u043e\u0432\u0430\u0442\u0435\u043b\u044f\u043c \u0441\u043e\u0446\u0441\u0435\u0442\u044c, \u0432 \u043b\u0435\u043d\u0442\u0435 \u043a\u043e\u0442\u043e\u0440\u043e\u0439 \u043f\u0440\u0435\u043e\u0431\u043b\u0430\u0434\u0430\u043b\u0438 \u0431\u044b \u044e\u043c\u043e\u0440\u0438\u0441\u0442\u0438\u0447\u0435\u0441\u043a\u0438\u0435 \u043f\u043e\u0441\u0442\u044b\"\n\n# 7. Locate the paragraph that holds the lone page break right before the\n# results table (it directly follows \"\u0421\u043f\u0430\u0441\u0438\u0431\u043e \u0437\u0430 \u043f\u043e\u043c\u043e\u0449\u044c!\").\n$paras = $d.Paragraphs\n$breakIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($t -like \"*$([char]12)*\") {\n        $prev = $paras.Item($i - 1).Range.Text\n        if ($prev -like \"*\u0421\u043f\u0430\u0441\u0438\u0431\u043e \u0437\u0430 \u043f\u043e\u043c\u043e\u0449\u044c*\") {\n            $breakIndex = $i\n            break\n        }\n    }\n}\nif ($breakIndex -eq -1) {\n    throw \"Could not locate the page-break paragraph before the table.\"\n}\n\n$breakPara = $paras.Item($breakIndex)\n\n# Insert the new (still empty) paragraph right after the page break first,\n# before the page-break paragraph mark itself is made bold/italic, so the\n# new paragraph does not inherit that formatting for its own mark.\n$breakPara.Range.InsertParagraphAfter()\n\n# The paragraph mark of the page-break paragraph picks up bold/italic\n# (matching the run that holds the page break).\n$paras = $d.Paragraphs\n$breakPara = $paras.Item($breakIndex)\n$breakPara.Range.Font.Bold = $true\n$breakPara.Range.Font.Italic = $true\n\n# Type the \"//TODO\" placeholder text into the new paragraph.\n$paras = $d.Paragraphs\n$todoPara = $paras.Item($breakIndex + 1)\n$todoPara.Range.Text = \"//TODO\"\n\n# Re-fetch a tight range spanning only the typed characters (not the\n# paragraph mark) so the bold/italic formatting lands on the run, not on\n# the paragraph's own mark.\n$paras = $d.Paragraphs\n$todoPara = $paras.Item($breakIndex + 1)\n$todoTextRange = $d.Range($todoPara.Range.Start, $todoPara.Range.Start + 6)\n$todoTextRange.Font.Bold = $true\n$todoTextRange.Font.Italic = $true\n$todoTextRange.LanguageID = \"en-US\"\n\n# Mark the whole paragraph (mark included) as English (US) as well.\n$paras = $d.Paragraphs\n$todoPara = $paras.Item($breakIndex + 1)\n$todoPara.Range.LanguageID = \"en-US\"\n"}
